$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).ColumnWidth = 15.42578125

$values = @(
    @(-0.07655054299873143, 0.07618794036417853),
    @(-0.030075889759855201, 0.028579679357184418),
    @(0.12679624710819937, -0.12744829742584685),
    @(-0.15654825632249469, 0.1557118895239924),
    @(-0.1497118896854257, 0.1480394558661926),
    @(-0.047159339197012873, 0.047132227089735856),
    @(-0.027132227290669775, 0.027091804774771688),
    @(-0.019542316813446092, 0.019514960558876027),
    @(-0.01351496072750713, 0.0134990793739167),
    @(-0.0074990795431233437, 0.007499153790355706),
    @(-0.0029991539559048874, 0.0029961207627273723),
    @(-0.0039552558384010617, 0.0039300925282770294),
    @(0.0020699073022072412, -0.0020829967881841327),
    @(0.014082996604198428, -0.014118973423443038),
    @(0.020118973254280803, -0.020177184032145412),
    @(0.026177183863575593, -0.026311219119081297),
    @(-0.0090043760414104668, 0.0089999998253640356),
    @(-0.036112050774026017, 0.036097027029295958),
    @(-0.027097027196643975, 0.027014086798309744),
    @(-0.018014086967184539, 0.018004324287954176),
    @(-0.0090043244570736647, 0.0089999998307099816),
    @(-0.14637592663335042, 0.145464375393761),
    @(-0.08463092724570398, 0.084126216311144475),
    @(-0.042126216564354557, 0.041999999745350181),
    @(-0.094880116835909689, 0.094638038039700945),
    @(-0.088638038207111691, 0.088327195976688699),
    @(-0.082327196145012493, 0.081268141329796872),
    @(-0.075268141501789287, 0.074534202095055235),
    @(-0.062534202284279417, 0.062169732103926734),
    @(-0.042169732313822372, 0.042020519429746717),
    @(-0.027020519629411766, 0.027001028082315415),
    @(-0.0060010282968185891, 0.0059999998216886397)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
